# Swap the contents of columns B:AC between each pair of rows listed below.
# Column A (the running index) is left untouched on both rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Saudi Arabia Division 1")

$pairs = @(
    @(38, 39),
    @(49, 50),
    @(65, 66),
    @(81, 82),
    @(125, 126),
    @(146, 147),
    @(149, 150),
    @(172, 173),
    @(191, 192),
    @(233, 234),
    @(238, 239),
    @(248, 249)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AC$r1")
    $range2 = $ws.Range("B$r2" + ":AC$r2")

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
